# ENW.xlsx update — adds three new test-case rows to the "Test Cases" sheet:
#   - row 26 (new): ENW024 / not-affiliated WOS->ENW alt-version case
#   - rows 31-32 (new, appended after the former last row): ENW025 / ENW026
#     (ENW icon from search results -> ENW sign-in / auto sign-in cases)
# Existing rows 26-29 shift down to 27-30 as a consequence of the row-26 insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Insert the new row 26 (pushes the former rows 26-29 down to 27-30)
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Insert()

# Values are entered in this precise order (A, C, B) so new shared-string
# entries land in the same sequence the workbook was originally authored in.
$ws.Range("A26").Value = "ENW024"
$ws.Range("C26").Value = "Verify that,user can navigate from WOS to ENW,if the user is not affiliated to a Customer in the market test group based on the WOS Customer Check, and should be sent to the alternate version of Endnote."
$ws.Range("B26").Value = "`nOPQA-3590"
$ws.Range("D26").Value = "Y"

# Formatting: match the look of the surrounding rows (thin-bordered cells,
# B/C wrap text, A/D/E do not).
$ws.Range("D11").Copy()
$ws.Range("A26").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("B26").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C26").PasteSpecial($xlPasteFormats)
$ws.Range("B2").Copy()
$ws.Range("D26").PasteSpecial($xlPasteFormats)
$ws.Range("D11").Copy()
$ws.Range("E26").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(26).RowHeight = 45

# ---------------------------------------------------------------------------
# 2) Append two new rows (31-32) after the (now shifted) last data row 30
# ---------------------------------------------------------------------------
$ws.Rows.Item("31:32").Insert()

# Row 31 values entered C, B, A; row 32 entered B, C, A - again matching the
# original authoring order so shared-string indices line up.
$ws.Range("C31").Value = "Verify that non personalized user who navigates from WOS to ENW by clicking the ENW icon from search results shall be sent to the ENW sign in screen for authentication"
$ws.Range("B31").Value = "QAOPQA-2733"
$ws.Range("A31").Value = "ENW025"
$ws.Range("D31").Value = "Y"

$ws.Range("B32").Value = "OPQA-2730"
$ws.Range("C32").Value = "Verify that personalized user who navigates from WOS to ENW by clicking the ENW icon from search results shall be automatically signed into ENW and taken to the full record"
$ws.Range("A32").Value = "ENW026"
$ws.Range("D32").Value = "Y"

foreach ($r in 31,32) {
    $ws.Range("A2").Copy()
    $ws.Range("A$r").PasteSpecial($xlPasteFormats)
    $ws.Range("D11").Copy()
    $ws.Range("B$r").PasteSpecial($xlPasteFormats)
    $ws.Range("A11").Copy()
    $ws.Range("C$r").PasteSpecial($xlPasteFormats)
    $ws.Range("B2").Copy()
    $ws.Range("D$r").PasteSpecial($xlPasteFormats)
    $ws.Range("D11").Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
    $ws.Rows.Item($r).RowHeight = 30
}

# ---------------------------------------------------------------------------
# 3) View state: selection moves to D29, and the sheet no longer needs to be
#    scrolled (topLeftCell gets cleared as a side effect of a fresh Select).
# ---------------------------------------------------------------------------
$ws.Range("D29").Select()

$wb.Save()
